# Rename the sheet from "new_departments_department" to "Sheet1".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# Reset the view's selection back to the default top-left cell (A1) instead
# of the previously-saved "A2:A6" selection.
[void]$ws.Range("A1").Select()
